# The commit removes the "Contd..." / "ISAM for developers and users can
# help them decide..." slide, which is slide #11 (SlideID 285) in the
# deck's slide order. Deleting it shifts every following slide up by one
# position, matching the updated <p:sldIdLst> in the target presentation.
$p = $ppt.ActivePresentation
$p.Slides.Item(11).Delete()
